# Release 2.2.0-ballot preparation edits for
# StructureDefinition-tddui-evaluation-reference.xlsx
#
# Changes applied:
#   1. Metadata!B3  (Version)          2.1.0 -> 2.2.0-ballot
#   2. Metadata!B8  (Date)             2025-12-18T17:25:31+00:00 -> 2025-12-19T08:32:44+00:00
#   3. Metadata!B18 (Base Definition)  append the |4.0.1 version to the canonical URL
#   4. Elements!K6  (Extension.value[x] Type(s)) append the |2.2.0-ballot version
#      to the Reference(...) canonical URL
#   5. Elements column K (11) widened to match the longer Type(s) text

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B3").Value = "2.2.0-ballot"
$wsMeta.Range("B8").Value = "2025-12-19T08:32:44+00:00"
$wsMeta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet --------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

$wsElem.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-questionnaire-response|2.2.0-ballot)
"

$wsElem.Columns.Item(11).ColumnWidth = 90.8
